$d = $word.ActiveDocument

$d.Content.Find.Execute("119÷9=13, 2", $false, $false, $false, $false, $false, $true, 1, $false, "240÷5=48, 0", 2)
$d.Content.Find.Execute("571÷5=114, 1", $false, $false, $false, $false, $false, $true, 1, $false, "532÷8=66, 4", 2)
$d.Content.Find.Execute("873÷4=218, 1", $false, $false, $false, $false, $false, $true, 1, $false, "826÷3=275, 1", 2)
$d.Content.Find.Execute("779÷3=259, 2", $false, $false, $false, $false, $false, $true, 1, $false, "722÷4=180, 2", 2)
$d.Content.Find.Execute("139÷6=23, 1", $false, $false, $false, $false, $false, $true, 1, $false, "546÷3=182, 0", 2)
$d.Content.Find.Execute("508÷7=72, 4", $false, $false, $false, $false, $false, $true, 1, $false, "522÷7=74, 4", 2)
$d.Content.Find.Execute("837÷7=119, 4", $false, $false, $false, $false, $false, $true, 1, $false, "497÷7=71, 0", 2)
$d.Content.Find.Execute("156÷8=19, 4", $false, $false, $false, $false, $false, $true, 1, $false, "156÷5=31, 1", 2)
$d.Content.Find.Execute("173÷2=86, 1", $false, $false, $false, $false, $false, $true, 1, $false, "659÷5=131, 4", 2)
$d.Content.Find.Execute("683÷3=227, 2", $false, $false, $false, $false, $false, $true, 1, $false, "191÷4=47, 3", 2)
$d.Content.Find.Execute("565÷5=113, 0", $false, $false, $false, $false, $false, $true, 1, $false, "196÷9=21, 7", 2)
$d.Content.Find.Execute("365÷6=60, 5", $false, $false, $false, $false, $false, $true, 1, $false, "363÷2=181, 1", 2)
$d.Content.Find.Execute("396÷8=49, 4", $false, $false, $false, $false, $false, $true, 1, $false, "357÷9=39, 6", 2)
$d.Content.Find.Execute("582÷7=83, 1", $false, $false, $false, $false, $false, $true, 1, $false, "377÷7=53, 6", 2)
$d.Content.Find.Execute("176÷6=29, 2", $false, $false, $false, $false, $false, $true, 1, $false, "962÷5=192, 2", 2)
$d.Content.Find.Execute("998÷2=499, 0", $false, $false, $false, $false, $false, $true, 1, $false, "709÷9=78, 7", 2)
$d.Content.Find.Execute("653÷2=326, 1", $false, $false, $false, $false, $false, $true, 1, $false, "644÷8=80, 4", 2)
$d.Content.Find.Execute("578÷7=82, 4", $false, $false, $false, $false, $false, $true, 1, $false, "312÷8=39, 0", 2)
$d.Content.Find.Execute("905÷2=452, 1", $false, $false, $false, $false, $false, $true, 1, $false, "525÷3=175, 0", 2)
$d.Content.Find.Execute("746÷9=82, 8", $false, $false, $false, $false, $false, $true, 1, $false, "770÷3=256, 2", 2)
$d.Content.Find.Execute("358÷6=59, 4", $false, $false, $false, $false, $false, $true, 1, $false, "766÷8=95, 6", 2)
$d.Content.Find.Execute("200÷4=50, 0", $false, $false, $false, $false, $false, $true, 1, $false, "374÷7=53, 3", 2)
$d.Content.Find.Execute("612÷8=76, 4", $false, $false, $false, $false, $false, $true, 1, $false, "592÷7=84, 4", 2)
$d.Content.Find.Execute("358÷8=44, 6", $false, $false, $false, $false, $false, $true, 1, $false, "492÷9=54, 6", 2)
$d.Content.Find.Execute("583÷7=83, 2", $false, $false, $false, $false, $false, $true, 1, $false, "321÷5=64, 1", 2)
